$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "text" column (B) -------------------------------------------------
# New wording replaces several of the survey-response strings (swapping "news
# story" -> "story" phrasing, adding new responses, dropping some old ones).
# The shared-strings table this runtime rebuilds on save keeps already-used
# strings in their existing relative order and appends brand-new ones in the
# order they are first assigned - so we deliberately touch the cells that
# introduce brand-new text first, in the exact order the final table needs
# them, before touching the cells that merely reuse still-live old text.

# 1) Cells that introduce the NEW strings, in the order they must appear:
$ws.Range("B2").Value  = "I agree with the story. "
$ws.Range("B6").Value  = "I disagree with the story. "
$ws.Range("B13").Value = "I see all of your points; I agree with the story too."
$ws.Range("B8").Value  = "I do agree with the story. "
$ws.Range("B5").Value  = "I have very often experienced such issues in the homes where I was."
$ws.Range("B7").Value  = "In the homes where I was, I have experienced none of the issues mentioned in the story."
$ws.Range("B10").Value = "I can warn people when they should check their Internet before connection shuts down. "
$ws.Range("B12").Value = "There have been many issues."
$ws.Range("B11").Value = "I have often warned people. "

# 2) Cells that just reuse text which already exists elsewhere in the sheet:
$ws.Range("B3").Value  = "I'm a robot connected to smart thermostats via the internet."
$ws.Range("B4").Value  = "I know a lot about the technology to evaluate its performance. "
$ws.Range("B9").Value  = "I can detect when Internet connectivity weakens. "
$ws.Range("B14").Value = "I have temperature sensors to detect when a room is too hot or too cold. "
$ws.Range("B15").Value = "I can fix it when a thermostat is not working correctly."

# 3) Brand-new row 16, appended at the bottom of the table.
$ws.Cells.Item(16, 1).Value = 3
$ws.Cells.Item(16, 2).Value = "I have always experienced temperature problems in the homes where I have been."

# --- Update the numeric columns (A and C) that moved around -----------------------
$ws.Range("C11").Value = 0.3
$ws.Range("A12").Value = 2
$ws.Range("C12").Value = 1.5
$ws.Range("C15").Value = 0.3

# --- View / selection tweaks --------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.Zoom = 140
$ws.Range("C17").Select()
